$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new appropriation entry (row 14), reusing row 13's formatting
# by copying it down (so the date/time number formats - and their
# underlying style records - are preserved exactly as the existing rows).
$ws.Rows("13:13").Copy() | Out-Null
$ws.Rows("14:14").Insert(-4121) | Out-Null

# New entry: 10/20/2013, 3h25 (0.1423611111111111 of a day)
$ws.Range("A14").Value = 41567
$ws.Range("B14").Value = 0.1423611111111111

# Reflect the new row in the current selection
$ws.Range("B4:B14").Select() | Out-Null
